# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the b.md file
# has been queued for a new handoff: status moves from "Handed back: in
# sync with en-US" to "Ready for handoff" (for both zh-cn and de-de),
# content duplicate flips to False, a new handoff file/datetime is
# recorded, and an error detail note about the stale handback is added.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (row 3 = b.md) ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Cells.Item(3, 5).Value = "Ready for handoff"            # E3 zh-cn status
$ov.Cells.Item(3, 6).Value = "Ready for handoff"            # F3 de-de status
$ov.Cells.Item(3, 7).Value = "2016-11-14 06:17:41"          # G3 latest HO xliff generate date

# ---- zh-cn sheet (row 3 = b.md) ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Cells.Item(3, 3).Value  = "Ready for handoff"                                              # C3 Status
# Leading apostrophe forces text (not boolean) storage for "False"; reset the
# style afterwards so we don't leave a stray quote-prefix format behind.
$zhF3 = $zh.Cells.Item(3, 6)
$zhF3.Value = "'False"                                                                         # F3 Content Duplicate
$zhF3.Style = "Normal"
$zh.Cells.Item(3, 7).Value  = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"            # G3 Latest Handoff File
$zh.Cells.Item(3, 8).Value  = "2016-11-14 06:17:27"                                            # H3 Latest Handoff Datetime
$zh.Cells.Item(3, 16).Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d256d7a009fa968f77dc2a92e2a32b31850becca/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/763b9d04475f8ea8daaad7b68d4598e04f170d8c/e2e/b.md."  # P3 Error Detail
# Column G is already a stored width of 40 in this sheet; reuse its COM
# ColumnWidth reading so column P lands on the same stored width=40 exactly
# (Excel's character-width<->pixel rounding makes a literal "=40" land on 40.83).
$zh.Columns.Item(16).ColumnWidth = $zh.Columns.Item(7).ColumnWidth()

# ---- de-de sheet (row 3 = b.md) ----
$de = $wb.Worksheets.Item("de-de")
$de.Cells.Item(3, 3).Value  = "Ready for handoff"                                              # C3 Status
$deF3 = $de.Cells.Item(3, 6)
$deF3.Value = "'False"                                                                         # F3 Content Duplicate
$deF3.Style = "Normal"
$de.Cells.Item(3, 7).Value  = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"            # G3 Latest Handoff File
$de.Cells.Item(3, 8).Value  = "2016-11-14 06:17:41"                                            # H3 Latest Handoff Datetime
$de.Cells.Item(3, 16).Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d256d7a009fa968f77dc2a92e2a32b31850becca/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/763b9d04475f8ea8daaad7b68d4598e04f170d8c/e2e/b.md."  # P3 Error Detail
$de.Columns.Item(16).ColumnWidth = $de.Columns.Item(7).ColumnWidth()
